$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the standalone "Investments" trial-balance row (old row 152) together
# with the January 2024 and February 2024 monthly cash-flow blocks (old rows
# 153-170). Everything below shifts up by 19 rows; Excel will also drop the
# now-unused "Investments" shared string automatically on save.
$ws.Rows("152:170").Delete() | Out-Null

# Restore the active selection to where the user was last working.
$ws.Range("H16").Select() | Out-Null

# A hidden sheet-scoped defined name marking the filter database range, as
# Excel leaves behind after toggling AutoFilter on the data range.
$nm = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$E`$286")
$nm.Visible = $false
